$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) "TextBox 243": "3.2 " + "upload patch"  ->  single run "3.2 upload patch"
# ---------------------------------------------------------------------------
$sh = $s.Shapes.Item("TextBox 243")
$tr = $sh.TextFrame.TextRange
$tr.Text = "TEMP"
$sh.TextFrame.TextRange.Text = "3.2 upload patch"

# ---------------------------------------------------------------------------
# 2) "TextBox 249": "5. " + "download patch"  ->  single run "5. download patch"
# ---------------------------------------------------------------------------
$sh = $s.Shapes.Item("TextBox 249")
$tr = $sh.TextFrame.TextRange
$tr.Text = "TEMP"
$sh.TextFrame.TextRange.Text = "5. download patch"

# ---------------------------------------------------------------------------
# 3) "TextBox 251": "6.1 test" <br> "6.2 " + "apply" -> "6.1 test" <br> "6.2 apply"
#    Only the last two runs (after the line break) are touched; the first run
#    and the line break must stay intact.
# ---------------------------------------------------------------------------
$sh = $s.Shapes.Item("TextBox 251")
$tr = $sh.TextFrame.TextRange
$tail = $tr.Characters(10, 9)
$tail.Text = "TEMP"
$sh.TextFrame.TextRange.Characters(10, 4).Text = "6.2 apply"

# ---------------------------------------------------------------------------
# 4) "TextBox 261": move/resize the box and split "3.1 pull" into
#    "3.1 " + "pull, merge"
# ---------------------------------------------------------------------------
$sh = $s.Shapes.Item("TextBox 261")
$sh.Left = 148.628428
$sh.Width = 90.716300
$tr = $sh.TextFrame.TextRange
$tr.Characters(5, 4).Text = "pull, merge"
